$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value updates derived from the diff (crypto price/volume refresh).
$updates = [ordered]@{
    'D2' = '67.513.20'
    'E2' = '  -0.28%  '
    'D3' = '3.778.17'
    'E3' = '  -0.36%  '
    'E4' = '  +0.05%  '
    'D5' = '597.40'
    'E5' = '  +0.31%  '
    'D6' = '164.58'
    'E6' = '  -1.43%  '
    'E7' = '  +0.00%  '
    'E8' = '  -1.17%  '
    'E9' = '  -1.03%  '
    'D10' = '0.449'
    'E10' = '  +0.21%  '
    'E11' = '  +1.33%  '
    'D12' = '0.0000248'
    'E12' = '  -2.10%  '
    'D13' = '35.57'
    'E13' = '  -1.36%  '
    'D14' = '4.414.21'
    'E14' = '  -0.27%  '
    'D15' = '3.787.61'
    'E15' = '  -0.46%  '
    'D16' = '67.553.08'
    'E16' = '  -0.12%  '
    'D17' = '18.25'
    'E17' = '  -0.85%  '
    'E18' = '  +1.63%  '
    'E19' = '  -0.95%  '
    'D20' = '460.18'
    'E20' = '  +0.23%  '
    'E21' = '  -2.23%  '
    'E22' = '  -0.34%  '
    'D23' = '0.0000146'
    'E23' = '  -5.44%  '
    'D24' = '82.35'
    'E24' = '  -1.12%  '
    'D25' = '11.96'
    'E25' = '  -0.74%  '
    'E26' = '  -1.40%  '
    'E27' = '  -0.43%  '
    'E28' = '  -0.53%  '
    'D29' = '3.925.60'
    'E29' = '  -0.35%  '
    'D30' = '7.43'
    'E30' = '  +2.55%  '
    'D31' = '2.65'
    'E31' = '  -4.51%  '
    'E32' = '  -2.89%  '
    'D33' = '28.96'
    'E33' = '  -2.31%  '
    'E34' = '  -0.33%  '
    'D35' = '8.96'
    'E35' = '  -1.33%  '
    'E36' = '  -1.52%  '
    'E37' = '  +0.28%  '
    'D38' = '0.989'
    'E38' = '  -0.42%  '
    'D39' = '3.23'
    'E39' = '  -4.38%  '
    'E40' = '  -0.80%  '
    'E41' = '  +0.11%  '
    'E42' = '  +0.03%  '
    'D43' = '43.57'
    'E43' = '  -1.27%  '
    'D44' = '47.46'
    'E44' = '  -1.24%  '
    'D45' = '0.296'
    'E45' = '  -0.77%  '
    'D46' = '151.40'
    'E46' = '  +1.22%  '
    'E47' = '  +0.35%  '
    'D48' = '1.37'
    'E48' = '  +8.27%  '
    'D49' = '27.16'
    'E49' = '  +1.33%  '
    'B50' = 'Bittensor'
    'C50' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D50' = '393.33'
    'E50' = '  +0.47%  '
    'B51' = 'Stacks'
    'C51' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D51' = '1.84'
    'E51' = '  +1.33%  '
}

# Every value in this sheet is stored as literal text (prices like "597.40"
# or "67.513.20", percentages like "  -0.28%  "). Excel normally auto-converts
# a plain decimal-looking string assigned via .Value into a real number, which
# would silently drop trailing zeros (e.g. "597.40" -> 597.4) or switch to
# scientific notation. Force the cell to Text format first whenever the new
# value would otherwise parse as a plain number, so it round-trips as a string.
foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $newValue = $updates[$ref]
    if ($newValue -match "^[+-]?\d+(\.\d+)?$") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $newValue
}
